$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# Row 52: only Taxonsorteringsordning (B) changes
$ws.Range("B52").Value = 90799

# Row 53
$ws.Range("A53").Value = 112111388
$ws.Range("B53").Value = 89104
$ws.Range("E53").Value = 5747
$ws.Range("F53").Value = "Läderdoftande fingersvamp"
$ws.Range("G53").Value = "Ramaria safraniolens"
$ws.Range("H53").Value = "Christian"

# Row 54
$ws.Range("A54").Value = 112111386
$ws.Range("B54").Value = 89094
$ws.Range("D54").Value = "VU"
$ws.Range("E54").Value = 256335
$ws.Range("F54").Value = "Taggfingersvamp"
$ws.Range("G54").Value = "Ramaria karstenii"
$ws.Range("H54").Value = "(Sacc. & P.Syd.) Corner"
$ws.Range("P54").Value = "Renkullmyren (Renkullmyren), Jmt"
$ws.Range("Q54").Value = 446734
$ws.Range("R54").Value = 7032709

# Row 55
$ws.Range("A55").Value = 112105381
$ws.Range("B55").Value = 89047
$ws.Range("D55").Value = "NT"
$ws.Range("E55").Value = 3286
$ws.Range("F55").Value = "Flattoppad klubbsvamp"
$ws.Range("G55").Value = "Clavariadelphus truncatus"
$ws.Range("H55").Value = "(Quél.) Donk"
$ws.Range("P55").Value = "Landverktjärnen (Landverktjärnen), Jmt"
$ws.Range("Q55").Value = 446564
$ws.Range("R55").Value = 7032716

# Row 56
$ws.Range("A56").Value = 112111378
$ws.Range("B56").Value = 83086
$ws.Range("E56").Value = 5589
$ws.Range("F56").Value = "Rödbrun klubbdyna"
$ws.Range("G56").Value = "Trichoderma nybergianum"
$ws.Range("H56").Value = "(T.Ulvinen & H.L.Chamb.) Jaklitsch & Voglmayr"
$ws.Range("Q56").Value = 446760
$ws.Range("R56").Value = 7032715

# Row 57
$ws.Range("A57").Value = 112111398
$ws.Range("B57").Value = 89114
$ws.Range("E57").Value = 5754
$ws.Range("F57").Value = "Gultoppig fingersvamp"
$ws.Range("G57").Value = "Ramaria testaceoflava"
$ws.Range("H57").Value = "(Bres.) Corner"
$ws.Range("Q57").Value = 446740
$ws.Range("R57").Value = 7032705

# Row 58
$ws.Range("A58").Value = 112213232
$ws.Range("B58").Value = 89104
$ws.Range("D58").Value = "VU"
$ws.Range("E58").Value = 5747
$ws.Range("F58").Value = "Läderdoftande fingersvamp"
$ws.Range("G58").Value = "Ramaria safraniolens"
$ws.Range("H58").Value = "Christian"
$ws.Range("Q58").Value = 446675
$ws.Range("R58").Value = 7032593

# Row 59
$ws.Range("A59").Value = 112213255
$ws.Range("B59").Value = 90480
$ws.Range("D59").Value = "LC"
$ws.Range("E59").Value = 4769
$ws.Range("F59").Value = "Svavelriska"
$ws.Range("G59").Value = "Lactarius scrobiculatus"
$ws.Range("H59").Value = "(Scop.:Fr.) Fr."
$ws.Range("Q59").Value = 446605
$ws.Range("R59").Value = 7032710

# Row 60: only Taxonsorteringsordning (B) changes
$ws.Range("B60").Value = 90835

# Row 61: only Taxonsorteringsordning (B) changes
$ws.Range("B61").Value = 88150
